# Applies updated "Price" (D) and "Volume(1h)" (E) figures from the
# Mon Jan 23 22:41:52 UTC 2023 GitHub Actions symbol-list refresh.
# Values are plain text in the source sheet (e.g. "306.06", "1.73%"),
# so each target cell is force-formatted as text before the write to
# stop Excel from auto-coercing the percent strings into numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Address -> new text value
$updates = [ordered]@{
    "D2" = "306.06"
    "E2" = "1.73%"
    "D3" = "36.40"
    "E3" = "-0.38%"
    "D4" = "5.065"
    "E4" = "1.52%"
    "D5" = "0.07922"
    "E5" = "3.16%"
    "D6" = "2.198"
    "E6" = "6.67%"
    "D7" = "8.019"
    "E7" = "1.27%"
    "D8" = "0.9308"
    "E8" = "1.37%"
    "D9" = "0.09865"
    "E9" = "2.16%"
    "E10" = "0.88%"
    "D11" = "0.09108"
    "E11" = "7.29%"
    "D12" = "0.03696"
    "E12" = "4.88%"
    "D13" = "0.09927"
    "E13" = "-0.27%"
    "D14" = "0.001446"
    "E14" = "-2.30%"
    "D15" = "0.005654"
    "E15" = "0.15%"
    "D16" = "3.466"
    "E16" = "0.04%"
    "D17" = "4.174"
    "E17" = "3.57%"
    "E18" = "8.75%"
    "D19" = "0.3371"
    "E19" = "-0.39%"
    "E20" = "1.48%"
    "D21" = "5.094"
    "E21" = "7.08%"
    "E22" = "-0.51%"
    "D23" = "0.04557"
    "E23" = "-0.80%"
    "D24" = "0.001239"
    "E24" = "0.65%"
    "E25" = "-5.91%"
    "D26" = "0.0001300"
    "E26" = "-7.09%"
    "D39" = "0.01944"
    "E39" = "10.73%"
    "D40" = "0.04930"
    "E40" = "7.11%"
    "D41" = "0.007796"
    "E41" = "4.32%"
    "D42" = "0.1398"
    "E42" = "0.68%"
    "D43" = "0.007810"
    "E43" = "1.07%"
    "D44" = "0.002112"
    "E44" = "-5.66%"
    "D45" = "0.01122"
    "E45" = "8.61%"
    "D46" = "0.00006232"
    "E46" = "-0.73%"
    "E47" = "-0.02%"
    "D48" = "52.23"
    "E48" = "48.45%"
    "D49" = "0.001801"
    "E49" = "-10.01%"
    "D50" = "0.00002101"
    "E50" = "-0.02%"
    "D51" = "0.0002001"
    "E51" = "-0.02%"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $updates[$addr]
}

# Drop the temporary text format so the cells fall back to the default
# (unstyled) cellXfs entry, matching the original, unstyled D/E cells.
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).ClearFormats()
}
